$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "Meta description: ..." paragraph that used to sit
#    right under the H1 title paragraph.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Cat Wilde and the Doom of Dead
#    Free | Slot Game Review" right before the final "Prompt: ..."
#    paragraph (i.e. right after the "Only 10 paylines ..." bullet).
# ------------------------------------------------------------------
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Only 10 paylines may limit betting options for some*") {
        $anchorIdx = $i
        break
    }
}

$anchorPara = $d.Paragraphs.Item($anchorIdx)
$anchorRange = $anchorPara.Range.Duplicate
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIdx + 1)
$newRange = $newPara.Range
$newRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cat Wilde and the Doom of Dead Free | Slot Game Review</w:t></w:r></w:p>')

# ------------------------------------------------------------------
# 3) Replace the long "Prompt: Create a cartoon-style..." italic text
#    with the new meta-description sentence.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptRange = $lastPara.Range
$oldText = "Prompt: Create a cartoon-style feature image for " + [char]34 + "Cat Wilde and the Doom of Dead" + [char]34 + " slot game with a happy Maya warrior wearing glasses as the main focus. The image should be visually stunning with bold colors and a playful, adventurous tone. In the center of the image, feature Cat Wilde, the adventurous explorer of ancient Egypt, standing confident and smiling with a fierce warrior outfit. Add some sand dunes and pyramids in the background, and have the Eye of Ra symbol, which acts as the wild in the game, hovering above Wilde's head and shining brightly. To Wilde's right, draw a happy Maya warrior wearing glasses and holding a compass and a sarcophagus. Make sure the warrior is looking pleased and excited about the treasure hunting adventure. In the top left corner of the image, add the game's title, " + [char]34 + "Cat Wilde and the Doom of Dead" + [char]34 + " in bold letters with the subtitle " + [char]34 + "Join Cat Wilde on an ancient Egyptian adventure" + [char]34 + " just below. Make sure the image is visually appealing, age-appropriate, and attention-grabbing enough to entice players to give the game a try."
$newText = "Explore ancient Egypt and win big in Cat Wilde and the Doom of Dead. Read our review and play for free today!"
$promptRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
